$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. Force text number format so numeric-looking
# strings (prices, hour codes) are preserved as text, matching the
# source data's original inline-string typing.
$updates = @(
    @{ Cell = "D2"; Value = "249.09" }
    @{ Cell = "G2"; Value = "5" }
    @{ Cell = "D3"; Value = "22.34" }
    @{ Cell = "G3"; Value = "5" }
    @{ Cell = "D4"; Value = "5.638" }
    @{ Cell = "G4"; Value = "5" }
    @{ Cell = "D5"; Value = "0.05600" }
    @{ Cell = "G5"; Value = "5" }
    @{ Cell = "D6"; Value = "3.390" }
    @{ Cell = "G6"; Value = "5" }
    @{ Cell = "D7"; Value = "6.488" }
    @{ Cell = "G7"; Value = "5" }
    @{ Cell = "G8"; Value = "5" }
    @{ Cell = "D9"; Value = "0.8015" }
    @{ Cell = "G9"; Value = "5" }
    @{ Cell = "G10"; Value = "5" }
    @{ Cell = "G11"; Value = "5" }
    @{ Cell = "D12"; Value = "0.07408" }
    @{ Cell = "G12"; Value = "5" }
    @{ Cell = "D13"; Value = "0.03224" }
    @{ Cell = "G13"; Value = "5" }
    @{ Cell = "B14"; Value = "ProBitToken" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" }
    @{ Cell = "D14"; Value = "0.1289" }
    @{ Cell = "E14"; Value = "13ProBitTokenPROB" }
    @{ Cell = "G14"; Value = "5" }
    @{ Cell = "B15"; Value = "BitrueCoin" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D15"; Value = "0.02994" }
    @{ Cell = "E15"; Value = "14BitrueCoinBTR" }
    @{ Cell = "G15"; Value = "5" }
    @{ Cell = "B16"; Value = "BitMartToken" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D16"; Value = "0.09264" }
    @{ Cell = "E16"; Value = "15BitMartTokenBMX" }
    @{ Cell = "G16"; Value = "5" }
    @{ Cell = "B17"; Value = "BitForexToken" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D17"; Value = "0.001668" }
    @{ Cell = "E17"; Value = "16BitForexTokenBF" }
    @{ Cell = "G17"; Value = "5" }
    @{ Cell = "B18"; Value = "MCDex" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D18"; Value = "3.249" }
    @{ Cell = "E18"; Value = "17MCDexMCB" }
    @{ Cell = "G18"; Value = "5" }
    @{ Cell = "B19"; Value = "CoinExToken" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" }
    @{ Cell = "D19"; Value = "0.04730" }
    @{ Cell = "E19"; Value = "18CoinExTokenCET" }
    @{ Cell = "G19"; Value = "5" }
    @{ Cell = "B20"; Value = "TigerCash" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D20"; Value = "0.006258" }
    @{ Cell = "E20"; Value = "19TigerCashTCH" }
    @{ Cell = "G20"; Value = "5" }
    @{ Cell = "B21"; Value = "BitKan" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan" }
    @{ Cell = "D21"; Value = "0.001047" }
    @{ Cell = "E21"; Value = "20BitKanKAN" }
    @{ Cell = "G21"; Value = "5" }
    @{ Cell = "B22"; Value = "HotbitToken" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" }
    @{ Cell = "D22"; Value = "0.003822" }
    @{ Cell = "E22"; Value = "21HotbitTokenHTB" }
    @{ Cell = "G22"; Value = "5" }
    @{ Cell = "B23"; Value = "NitroEx" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx" }
    @{ Cell = "D23"; Value = "0.0001499" }
    @{ Cell = "E23"; Value = "22NitroExNTX" }
    @{ Cell = "G23"; Value = "5" }
    @{ Cell = "B24"; Value = "UpBots" }
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt" }
    @{ Cell = "D24"; Value = "0.0004774" }
    @{ Cell = "E24"; Value = "23UpBotsUBXT" }
    @{ Cell = "G24"; Value = "5" }
    @{ Cell = "B25"; Value = "LEO" }
    @{ Cell = "C25"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D25"; Value = "3.979" }
    @{ Cell = "E25"; Value = "24LEOLEO" }
    @{ Cell = "G25"; Value = "5" }
    @{ Cell = "B26"; Value = "BTSEToken" }
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "D26"; Value = "2.119" }
    @{ Cell = "E26"; Value = "25BTSETokenBTSE" }
    @{ Cell = "G26"; Value = "5" }
    @{ Cell = "B27"; Value = "BitpandaEcosystemToken" }
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Cell = "D27"; Value = "0.3311" }
    @{ Cell = "E27"; Value = "26BitpandaEcosystemTokenBEST" }
    @{ Cell = "G27"; Value = "5" }
    @{ Cell = "G28"; Value = "5" }
    @{ Cell = "G29"; Value = "5" }
    @{ Cell = "G30"; Value = "5" }
    @{ Cell = "G31"; Value = "5" }
    @{ Cell = "G32"; Value = "5" }
    @{ Cell = "G33"; Value = "5" }
    @{ Cell = "G34"; Value = "5" }
    @{ Cell = "G35"; Value = "5" }
    @{ Cell = "G36"; Value = "5" }
    @{ Cell = "G37"; Value = "5" }
    @{ Cell = "G38"; Value = "5" }
    @{ Cell = "G39"; Value = "5" }
    @{ Cell = "D40"; Value = "0.04199" }
    @{ Cell = "G40"; Value = "5" }
    @{ Cell = "D41"; Value = "0.007005" }
    @{ Cell = "G41"; Value = "5" }
    @{ Cell = "D42"; Value = "0.1047" }
    @{ Cell = "G42"; Value = "5" }
    @{ Cell = "D43"; Value = "0.003090" }
    @{ Cell = "G43"; Value = "5" }
    @{ Cell = "D44"; Value = "0.009848" }
    @{ Cell = "G44"; Value = "5" }
    @{ Cell = "D45"; Value = "0.00005668" }
    @{ Cell = "G45"; Value = "5" }
    @{ Cell = "G46"; Value = "5" }
    @{ Cell = "D47"; Value = "0.6795" }
    @{ Cell = "G47"; Value = "5" }
    @{ Cell = "D48"; Value = "0.02975" }
    @{ Cell = "G48"; Value = "5" }
    @{ Cell = "G49"; Value = "5" }
    @{ Cell = "G50"; Value = "5" }
    @{ Cell = "G51"; Value = "5" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
